$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44270
$ws.Range("M2").Value = 85
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("Q2").Value = '$/caja 14 kilos granel'
$ws.Range("S2").Value = 857

$ws.Range("D3").Value = 44245
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/caja 15 kilos granel'
$ws.Range("T3").Value = 15

$ws.Range("D4").Value = 44320
$ws.Range("M4").Value = 45

$ws.Range("D5").Value = 44271
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("R5").Value = 'Provincia del Elquí'
$ws.Range("S5").Value = 857

$ws.Range("D6").Value = 44239
$ws.Range("M6").Value = 70

$ws.Range("D7").Value = 44252
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = '$/caja 14 kilos empedrada'
$ws.Range("T7").Value = 14

$ws.Range("D8").Value = 44238
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 15 kilos granel'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 15

$ws.Range("D9").Value = 44313
$ws.Range("M9").Value = 36

$ws.Range("D11").Value = 44322
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 1000

$ws.Range("D13").Value = 44316
$ws.Range("M13").Value = 48

$ws.Range("D14").Value = 44315
$ws.Range("M14").Value = 65

$ws.Range("D15").Value = 44278
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 13000
$ws.Range("Q15").Value = '$/caja 14 kilos empedrada'
$ws.Range("R15").Value = 'Provincia del Elquí'
$ws.Range("S15").Value = 929

$ws.Range("D16").Value = 44312
$ws.Range("M16").Value = 68
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("Q16").Value = '$/caja 14 kilos granel'
$ws.Range("T16").Value = 14

$ws.Range("D17").Value = 44314
$ws.Range("M17").Value = 56
$ws.Range("N17").Value = 14000
$ws.Range("O17").Value = 14000
$ws.Range("P17").Value = 14000
$ws.Range("Q17").Value = '$/caja 14 kilos granel'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 1000

$ws.Range("D18").Value = 44260
$ws.Range("M18").Value = 56
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 13000
$ws.Range("Q18").Value = '$/caja 14 kilos empedrada'
$ws.Range("R18").Value = 'Provincia del Elquí'
$ws.Range("S18").Value = 929

$ws.Range("D19").Value = 44323
$ws.Range("M19").Value = 60

Write-Output "Edits applied"
